# Review [test]solar_wind ArcGIS projects (VEN, GHA, GEO, FRA, USA, CHN)
#
# Inserts a new data row for Venezuela into the country table on Sheet1,
# between Uzbekistan (row 184) and Virgin Islands (U.S.) (previously row
# 185). All subsequent rows shift down by one (185 -> 186, ..., 193 -> 194).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push row 185 (and everything after it) down by one row, creating a blank
# row 185 for the new Venezuela entry.
$ws.Rows("185:185").Insert()

# Identifying columns.
$ws.Range("A185").Value = "Venezuela"
$ws.Range("B185").Value = "VEN"
$ws.Range("C185").Value = "Latin America & Caribbean"
$ws.Range("D185").Value = "Upper middle income"

# Population figures.
$ws.Range("E185").Value = 28405.543
$ws.Range("F185").Value = 29151.81
$ws.Range("G185").Value = 31094.883
$ws.Range("H185").Value = 1.02627187940044
$ws.Range("I185").Value = 1.094676591818716

# 2024 generation capacity (GW) and generation (MWh).
$ws.Range("J185").Value = 16910
$ws.Range("K185").Value = 1.1
$ws.Range("L185").Value = 50
$ws.Range("M185").Value = 60
$ws.Range("N185").Value = 0
$ws.Range("O185").Value = 17250
$ws.Range("P185").Value = 64980000.00000001
$ws.Range("Q185").Value = 2015.22013622223
$ws.Range("R185").Value = 20000
$ws.Range("S185").Value = 296443.5448480757
$ws.Range("T185").Value = 0
$ws.Range("U185").Value = 17930000

# Renewables-share target.
$ws.Range("W185").Value = 62.82165000000001

# 2030 generation (MWh).
$ws.Range("AC185").Value = 64980000.00000001
$ws.Range("AD185").Value = 2015.22013622223
$ws.Range("AE185").Value = 20000
$ws.Range("AF185").Value = 296443.5448480757
$ws.Range("AG185").Value = 0
$ws.Range("AH185").Value = 16276119.88780221

# 2050 generation (MWh).
$ws.Range("AI185").Value = 69445887.83033092
$ws.Range("AJ185").Value = 2153.720399099926
$ws.Range("AK185").Value = 21374.54226849213
$ws.Range("AL185").Value = 316817.2539788418
$ws.Range("AM185").Value = 0
$ws.Range("AN185").Value = 11788345.30580914
